$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sequence")
# "Sequence" is also the active sheet in this workbook ($wb.ActiveSheet)

$ws.Range("B18").Value = "optimiser"
$ws.Range("A19").Value = "Nesterov Accelerated Gradient (NAG)"
$ws.Range("B19").Value = "optimiser"
$ws.Range("A20").Value = "AdaGrad"
$ws.Range("B20").Value = "optimiser"

[void]$ws.Range("D21").Select()
